# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q3" and before "总计",
#    populated with the per-fund holdings for the quarter.
# 2. Insert a new summary row at the top of the "总计" sheet's data for
#    "2022-Q1" (pushing the existing "2021-Q3" summary row down).

$wb = $excel.ActiveWorkbook

$wsQ3 = $wb.Worksheets("2021-Q3")

# --- 1. Create the new "2022-Q1" sheet right after "2021-Q3" -------------
$wsQ1 = $wb.Worksheets.Add($null, $wsQ3)
$wsQ1.Name = "2022-Q1"

# Re-fetch every sheet reference after inserting a sheet: handles grabbed
# before the insert can end up pointing at the wrong (shifted) tab.
$wsQ3 = $wb.Worksheets("2021-Q3")
$wsTotal = $wb.Worksheets("总计")

# Helper: write a column of values as genuine text (no numeric coercion,
# no left-over number-format) by staging them in a scratch column, copying
# only the *values*, then deleting the scratch column again.
function Set-TextColumn {
    param($ws, [string]$destCol, [int]$rowStart, [object[]]$values)

    $n = $values.Length
    $rowEnd = $rowStart + $n - 1
    $helper = $ws.Range("ZZ$rowStart" + ":ZZ$rowEnd")
    for ($i = 0; $i -lt $n; $i++) {
        $r = $rowStart + $i
        $ws.Range("ZZ$r").Value = "'" + $values[$i]
    }
    $helper.Copy()
    $ws.Range($destCol + "$rowStart" + ":" + $destCol + "$rowEnd").PasteSpecial(-4163) # xlPasteValues
    $helper.ClearContents()
}

# Header row text.
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Copy the header cell formatting (bold, centered, bordered) from the
# "总计" sheet's header row onto the new header row.
foreach ($col in @("B1", "C1", "D1", "E1", "F1", "G1", "H1")) {
    $wsTotal.Range("B1").Copy()
    $wsQ1.Range($col).PasteSpecial(-4122) # xlPasteFormats
}

# Row-index column (numeric) and rank column (numeric).
$wsQ1.Range("A2").Value = 0
$wsQ1.Range("A3").Value = 1
$wsQ1.Range("A4").Value = 2
$wsQ1.Range("H2").Value = 10
$wsQ1.Range("H3").Value = 10
$wsQ1.Range("H4").Value = 3

# Fund code / name / scale / position columns (kept as text, matching the
# other sheets' data rows).
Set-TextColumn $wsQ1 "B" 2 @("014746", "014747", "005269")
Set-TextColumn $wsQ1 "D" 2 @("5.05", "2.23", "0.33")
Set-TextColumn $wsQ1 "E" 2 @("53.79", "53.79", "37.77")
Set-TextColumn $wsQ1 "F" 2 @("1.92", "1.92", "1.22")
Set-TextColumn $wsQ1 "G" 2 @("0.0970", "0.0428", "0.0040")

$wsQ1.Range("C2").Value = "贝莱德港股通远景视野混合A"
$wsQ1.Range("C3").Value = "贝莱德港股通远景视野混合C"
$wsQ1.Range("C4").Value = "华泰柏瑞港股通量化灵活配置混合"

# Drop the scratch column entirely so no trace (values/format/dimension) of
# it survives.
$wsQ1.Range("ZZ1").EntireColumn.Delete()

# Copy the row-number column formatting (bold, centered, bordered) from
# "总计"'s A-column onto the new sheet's A2:A4.
foreach ($row in @("A2", "A3", "A4")) {
    $wsTotal.Range("A2").Copy()
    $wsQ1.Range($row).PasteSpecial(-4122) # xlPasteFormats
}

# --- 2. Insert the "2022-Q1" summary row into "总计" ----------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.14

# The previously-existing summary row shifted from row 2 to row 3; restore
# its row-index counter (it keeps its old "0" value after the shift).
$wsTotal.Range("A3").Value = 1

# Match the bold/centered/bordered style used on the other "A" summary cell.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122) # xlPasteFormats

Write-Output "2022-Q1 sheet added; summary sheet updated"
